# (2017/12/16) ESTRUTURA DAS APOSTAS - POT ODDS
# Adds the TURN / RIVER betting rounds (rows 5-10) under the existing
# FLOP row, along with the pot / pot-odds formulas for each player.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - FLOP (A5 already held "FLOP" but the row had no bet data yet)
$ws.Range("A5").Value = "FLOP"
$ws.Range("B5").Value = 30
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Formula = "=SUM(B5:F5)+G4"
$ws.Range("H5").Formula = "=(B5*100)/(B5+G5)"
$ws.Range("I5").Formula = "=(C5*100)/(C5+G5)"
$ws.Range("J5").Formula = "=(D5*100)/(D5+G5)"
$ws.Range("K5").Formula = "=(E5*100)/(E5+G5)"
$ws.Range("L5").Formula = "=(F5*100)/(F5+G5)"

# Row 6 - continuation of FLOP betting (no label in column A)
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = 60
$ws.Range("G6").Formula = "=SUM(B6:F6)+G5"
$ws.Range("H6").Formula = "=(B6*100)/(B6+G6)"
$ws.Range("I6").Formula = "=(C6*100)/(C6+G6)"
$ws.Range("J6").Formula = "=(D6*100)/(D6+G6)"
$ws.Range("K6").Formula = "=(E6*100)/(E6+G6)"
$ws.Range("L6").Formula = "=(F6*100)/(F6+G6)"

# Row 7 - TURN
$ws.Range("A7").Value = "TURN"
$ws.Range("B7").Value = 200
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 150
$ws.Range("G7").Formula = "=SUM(B7:F7)+G6"
$ws.Range("H7").Formula = "=(B7*100)/(B7+G7)"
$ws.Range("I7").Formula = "=(C7*100)/(C7+G7)"
$ws.Range("J7").Formula = "=(D7*100)/(D7+G7)"
$ws.Range("K7").Formula = "=(E7*100)/(E7+G7)"
$ws.Range("L7").Formula = "=(F7*100)/(F7+G7)"

# Row 8 - continuation of TURN betting (no label in column A)
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = 50
$ws.Range("G8").Formula = "=SUM(B8:F8)+G7"
$ws.Range("H8").Formula = "=(B8*100)/(B8+G8)"
$ws.Range("I8").Formula = "=(C8*100)/(C8+G8)"
$ws.Range("J8").Formula = "=(D8*100)/(D8+G8)"
$ws.Range("K8").Formula = "=(E8*100)/(E8+G8)"
$ws.Range("L8").Formula = "=(F8*100)/(F8+G8)"

# Row 9 - RIVER
$ws.Range("A9").Value = "RIVER"
$ws.Range("B9").Value = 200
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 50
$ws.Range("F9").Value = 150
$ws.Range("G9").Formula = "=SUM(B9:F9)+G8"
$ws.Range("H9").Formula = "=(B9*100)/(B9+G9)"
$ws.Range("I9").Formula = "=(C9*100)/(C9+G9)"
$ws.Range("J9").Formula = "=(D9*100)/(D9+G9)"
$ws.Range("K9").Formula = "=(E9*100)/(E9+G9)"
$ws.Range("L9").Formula = "=(F9*100)/(F9+G9)"

# Row 10 - continuation of RIVER betting (no label in column A)
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 50
$ws.Range("G10").Formula = "=SUM(B10:F10)+G9"
$ws.Range("H10").Formula = "=(B10*100)/(B10+G10)"
$ws.Range("I10").Formula = "=(C10*100)/(C10+G10)"
$ws.Range("J10").Formula = "=(D10*100)/(D10+G10)"
$ws.Range("K10").Formula = "=(E10*100)/(E10+G10)"
$ws.Range("L10").Formula = "=(F10*100)/(F10+G10)"

# Match the author's final cursor position
$ws.Range("L10").Select()
